$d = $word.ActiveDocument
$s = $d.Shapes.Item(3)
$r = $s.TextFrame.TextRange
$paras = $r.Paragraphs
Write-Output ("count: " + $paras.Count)
for ($i = 1; $i -le 6; $i++) {
  $p = $paras.Item($i)
  $t = $p.Range.Text
  Write-Output ("p" + $i + " len=" + $t.Length + " text=[" + $t + "]")
}
